# ------------------------------------------------------------------
# Applies the "ADDITIONAL SCRAPING" edit:
#  - New sheet "Player Info" (first position)
#  - Existing "ODI Batting" sheet: rename D column MATCH_CARD_LINK -> MATCH_CODE
#    (values become bare match codes instead of full URLs), and drop the
#    placeholder-empty B (INNING_NUMBER-less) cells that used to be written out.
#  - Existing "ODI Bowling" sheet: rename B column MATCH_CARD_LINK -> MATCH_CODE
#    (values become bare match codes instead of full URLs).
#  - New sheet "ODI Batting Extra" (last position) with additional per-match
#    batting detail columns.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$urlPrefix = "http://www.howstat.com/cricket/Statistics/Matches/MatchScorecard_ODI.asp?MatchCode="

# --------------------------------------------------------------
# 1. Locate the two pre-existing sheets (by their current names)
# --------------------------------------------------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# --------------------------------------------------------------
# 2. Update "ODI Batting": D1 header + D2:D100 values, drop blank B cells
# --------------------------------------------------------------
$battingWs.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingRange = $battingWs.Range("D2:D100")
$battingRange.NumberFormat = "@"
for ($r = 2; $r -le 100; $r++) {
    $cell = $battingWs.Cells.Item($r, 4)
    $orig = $cell.Value2
    $code = $orig.Replace($urlPrefix, "")
    $cell.Value = $code
}

$rowsToClearB = @(5,6,7,8,12,14,15,16,17,18,19,20,22,23,24,25,28,29,31,32,36,41,45,46,48,52,53,55,56,57,58,59,60,62,63,65,66,70,71,72,76,78,79,80,81,82,83,84,85,86,92,93,94,97)
foreach ($r in $rowsToClearB) {
    $battingWs.Cells.Item($r, 2).ClearContents()
}

# --------------------------------------------------------------
# 3. Update "ODI Bowling": B1 header + B2:B100 values
# --------------------------------------------------------------
$bowlingWs.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingRange = $bowlingWs.Range("B2:B100")
$bowlingRange.NumberFormat = "@"
for ($r = 2; $r -le 100; $r++) {
    $cell = $bowlingWs.Cells.Item($r, 2)
    $orig = $cell.Value2
    $code = $orig.Replace($urlPrefix, "")
    $cell.Value = $code
}

# --------------------------------------------------------------
# 4. Add new "Player Info" sheet as the first sheet
# --------------------------------------------------------------
$playerInfoWs = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$playerInfoWs.Name = "Player Info"

$piHeader = $playerInfoWs.Range("A1:D1")
$piHeader.Value = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
$battingWs.Range("A1:D1").Copy() | Out-Null
$piHeader.PasteSpecial(-4122) | Out-Null

$playerInfoWs.Range("A2:D2").NumberFormat = "@"
$playerInfoWs.Cells.Item(2, 1).Value = "3924"
$playerInfoWs.Cells.Item(2, 2).Value = "Trent Alexander Boult"
$playerInfoWs.Cells.Item(2, 3).Value = "Right Handed"
$playerInfoWs.Cells.Item(2, 4).Value = "Left Arm Fast Medium"

# --------------------------------------------------------------
# 5. Add new "ODI Batting Extra" sheet as the last sheet
# --------------------------------------------------------------
$lastIdx = $wb.Worksheets.Count
$extraWs = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIdx))
$extraWs.Name = "ODI Batting Extra"

$exHeader = $extraWs.Range("A1:F1")
$exHeader.Value = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
$battingWs.Range("A1:D1").Copy() | Out-Null
$exHeader.PasteSpecial(-4122) | Out-Null

$extraWs.Range("A2:A21").NumberFormat = "@"
$extraWs.Range("C2:F21").NumberFormat = "@"

$extraMatchCodes = @(4305,4311,4315,4328,4333,4337,4341,4346,4353,4355,4423,4452,4453,4455,4636,4639,4642,4647,4648,4649)
$extraManOfMatch = @("NO","NO","NO","NO","NO","NO","NO","NO","NO","NO","NO","","","","","","","","","")

for ($i = 0; $i -lt $extraMatchCodes.Length; $i++) {
    $r = $i + 2
    $extraWs.Cells.Item($r, 1).Value = [string]$extraMatchCodes[$i]
    $extraWs.Cells.Item($r, 2).Value = ""
    $extraWs.Cells.Item($r, 3).Value = ""
    $extraWs.Cells.Item($r, 4).Value = ""
    $extraWs.Cells.Item($r, 5).Value = ""
    $extraWs.Cells.Item($r, 6).Value = $extraManOfMatch[$i]
}

# Batting position / boundary-count / percentage detail known for a few matches
$extraWs.Cells.Item(9, 2).Value = 11   # 4341 -> BATTING_POSITION 11
$extraWs.Cells.Item(10, 2).Value = 11  # 4346 -> BATTING_POSITION 11
$extraWs.Cells.Item(10, 3).Value = "0"
$extraWs.Cells.Item(10, 4).Value = "0"
$extraWs.Cells.Item(10, 5).Value = "0.53%"
$extraWs.Cells.Item(11, 2).Value = 11  # 4353 -> BATTING_POSITION 11
$extraWs.Cells.Item(11, 3).Value = "2"
$extraWs.Cells.Item(11, 4).Value = "0"
$extraWs.Cells.Item(11, 5).Value = "7.55%"

# --------------------------------------------------------------
# 6. Final sheet order: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# --------------------------------------------------------------
$playerInfoFinal = $wb.Worksheets.Item("Player Info")
$playerInfoFinal.Move($wb.Worksheets.Item(1))

$battingFinal = $wb.Worksheets.Item("ODI Batting")
$bowlingFinal = $wb.Worksheets.Item("ODI Bowling")
$bowlingFinal.Move($null, $battingFinal)

$extraFinal = $wb.Worksheets.Item("ODI Batting Extra")
$extraFinal.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

Write-Output "Sheets now:"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Output $wb.Worksheets.Item($i).Name
}
